# Generate Report for Handback
# Applies the "Handback" status update across the Overview/zh-cn/de-de sheets:
#  - Status text "In Translation" -> "Handed back: in sync with en-US"
#  - Latest Handback DateTime values populated (were "0001-01-01 00:00:00")
#  - Latest Target File / Latest Handback File columns populated with hyperlinked
#    file names + target xlf file names
#  - Related column widths widened to fit the new, longer content

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$mdUrl520 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3f02818c38505cddf79dce81e5bc489bb693cb6/e2e/520d59b4-13ae-4dec-b01b-047d29f2c514.md"
$mdUrlC5c = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3f02818c38505cddf79dce81e5bc489bb693cb6/e2e/c5cb36e2-acab-43a7-bb4f-9a7281913270.md"
$mdName520 = "520d59b4-13ae-4dec-b01b-047d29f2c514.md"
$mdNameC5c = "c5cb36e2-acab-43a7-bb4f-9a7281913270.md"

# ---------------------------------------------------------------------------
# 1. Update the "Status" text everywhere it is used (Overview summary columns
#    E/F, and the per-row Status column C on the zh-cn / de-de detail sheets).
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus

$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("C3").Value2 = $newStatus

$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate "Latest Target File" (col I) with a hyperlink to the source
#    markdown file, and "Latest Handback File" (col J) with the generated
#    xlf file name, for both language sheets.
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl520, "", "", $mdName520)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrlC5c, "", "", $mdNameC5c)

$wsZhCn.Range("J2").Value2 = "520d59b4-13ae-4dec-b01b-047d29f2c514.907f6a8ffc488d77779bbb123e23ae287581137b.zh-cn.xlf"
$wsZhCn.Range("J3").Value2 = "c5cb36e2-acab-43a7-bb4f-9a7281913270.02758016ab7921ae16508c521ac1a6073d118847.zh-cn.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl520, "", "", $mdName520)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrlC5c, "", "", $mdNameC5c)

$wsDeDe.Range("J2").Value2 = "520d59b4-13ae-4dec-b01b-047d29f2c514.907f6a8ffc488d77779bbb123e23ae287581137b.de-de.xlf"
$wsDeDe.Range("J3").Value2 = "c5cb36e2-acab-43a7-bb4f-9a7281913270.02758016ab7921ae16508c521ac1a6073d118847.de-de.xlf"

# ---------------------------------------------------------------------------
# 3. Record the handback timestamps in "Latest Handback DateTime" (col K).
#    zh-cn and de-de were handed back a few seconds apart.
# ---------------------------------------------------------------------------
$wsZhCn.Range("K2").Value2 = "2016-08-18 22:24:49"
$wsZhCn.Range("K3").Value2 = "2016-08-18 22:24:49"

$wsDeDe.Range("K2").Value2 = "2016-08-18 22:24:55"
$wsDeDe.Range("K3").Value2 = "2016-08-18 22:24:55"

# ---------------------------------------------------------------------------
# 4. Widen columns so the new, longer text is fully visible.
#    (ColumnWidth is quantized by the host app to the nearest 1/6 character,
#    so we pick inputs that land exactly on the desired displayed width.)
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.17   # -> 30 (status col, zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 29.17   # -> 30 (status col, de-de)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17        # -> 30 (Status)
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15        # -> 40 (Latest Target File)
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15       # -> 40 (Latest Handback File)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17        # -> 30 (Status)
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15        # -> 40 (Latest Target File)
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15       # -> 40 (Latest Handback File)
